# Fruta / hortaliza, semanal
# Apply updated weekly price data for "Hortaliza, Feria Lagunitas de Puerto Montt - Alcachofa"
# Columns: D=Fecha, H=Variedad, I=Calidad, J=Volumen, K=Precio minimo, L=Precio maximo,
#          M=Precio promedio ponderado, N=Unidad de comercializacion, P=Precio $/Kg, Q=Kg o Unidades

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = 44446
$ws.Cells.Item(2, 11).Value = 16000
$ws.Cells.Item(2, 12).Value = 16000
$ws.Cells.Item(2, 13).Value = 16000
$ws.Cells.Item(2, 16).Value = 400

# Row 3
$ws.Cells.Item(3, 4).Value = 44161
$ws.Cells.Item(3, 10).Value = 30
$ws.Cells.Item(3, 14).Value = "$/caja 40 unidades"
$ws.Cells.Item(3, 16).Value = 275
$ws.Cells.Item(3, 17).Value = 40

# Row 4
$ws.Cells.Item(4, 4).Value = 44162
$ws.Cells.Item(4, 8).Value = "Madrigal"
$ws.Cells.Item(4, 10).Value = 50
$ws.Cells.Item(4, 11).Value = 10000
$ws.Cells.Item(4, 12).Value = 10000
$ws.Cells.Item(4, 13).Value = 10000
$ws.Cells.Item(4, 14).Value = "$/caja 40 unidades"
$ws.Cells.Item(4, 16).Value = 250
$ws.Cells.Item(4, 17).Value = 40

# Row 5
$ws.Cells.Item(5, 4).Value = 44488
$ws.Cells.Item(5, 8).Value = "Madrigal"
$ws.Cells.Item(5, 10).Value = 120
$ws.Cells.Item(5, 11).Value = 12000
$ws.Cells.Item(5, 12).Value = 12000
$ws.Cells.Item(5, 13).Value = 12000
$ws.Cells.Item(5, 14).Value = "$/caja 40 unidades"
$ws.Cells.Item(5, 16).Value = 300
$ws.Cells.Item(5, 17).Value = 40

# Row 6
$ws.Cells.Item(6, 4).Value = 44481
$ws.Cells.Item(6, 9).Value = "Segunda"
$ws.Cells.Item(6, 10).Value = 120
$ws.Cells.Item(6, 11).Value = 11000
$ws.Cells.Item(6, 12).Value = 11000
$ws.Cells.Item(6, 13).Value = 11000
$ws.Cells.Item(6, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(6, 16).Value = 220
$ws.Cells.Item(6, 17).Value = 50

# Row 7
$ws.Cells.Item(7, 4).Value = 44421
$ws.Cells.Item(7, 11).Value = 16500
$ws.Cells.Item(7, 12).Value = 16500
$ws.Cells.Item(7, 13).Value = 16500
$ws.Cells.Item(7, 16).Value = 550

# Row 8
$ws.Cells.Item(8, 4).Value = 44495
$ws.Cells.Item(8, 10).Value = 130
$ws.Cells.Item(8, 11).Value = 11000
$ws.Cells.Item(8, 12).Value = 11000
$ws.Cells.Item(8, 13).Value = 11000
$ws.Cells.Item(8, 16).Value = 275

# Row 9
$ws.Cells.Item(9, 4).Value = 44176
$ws.Cells.Item(9, 9).Value = "Primera"
$ws.Cells.Item(9, 10).Value = 80
$ws.Cells.Item(9, 14).Value = "$/caja 40 unidades"
$ws.Cells.Item(9, 16).Value = 275
$ws.Cells.Item(9, 17).Value = 40

# Row 10
$ws.Cells.Item(10, 4).Value = 44400
$ws.Cells.Item(10, 10).Value = 70
$ws.Cells.Item(10, 11).Value = 15000
$ws.Cells.Item(10, 12).Value = 15000
$ws.Cells.Item(10, 13).Value = 15000
$ws.Cells.Item(10, 16).Value = 500

# Row 11
$ws.Cells.Item(11, 4).Value = 44484
$ws.Cells.Item(11, 10).Value = 110
$ws.Cells.Item(11, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(11, 16).Value = 220
$ws.Cells.Item(11, 17).Value = 50

# Row 13
$ws.Cells.Item(13, 4).Value = 44390
$ws.Cells.Item(13, 8).Value = "Española"
$ws.Cells.Item(13, 10).Value = 80
$ws.Cells.Item(13, 11).Value = 16000
$ws.Cells.Item(13, 12).Value = 16000
$ws.Cells.Item(13, 13).Value = 16000
$ws.Cells.Item(13, 14).Value = "$/caja 30 unidades"
$ws.Cells.Item(13, 16).Value = 533
$ws.Cells.Item(13, 17).Value = 30

# Row 14
$ws.Cells.Item(14, 4).Value = 44386
$ws.Cells.Item(14, 8).Value = "Española"
$ws.Cells.Item(14, 10).Value = 30
$ws.Cells.Item(14, 11).Value = 15000
$ws.Cells.Item(14, 12).Value = 15000
$ws.Cells.Item(14, 13).Value = 15000
$ws.Cells.Item(14, 14).Value = "$/caja 30 unidades"
$ws.Cells.Item(14, 16).Value = 500
$ws.Cells.Item(14, 17).Value = 30

# Row 15
$ws.Cells.Item(15, 4).Value = 44166
$ws.Cells.Item(15, 10).Value = 80
$ws.Cells.Item(15, 11).Value = 10000
$ws.Cells.Item(15, 12).Value = 10000
$ws.Cells.Item(15, 13).Value = 10000
$ws.Cells.Item(15, 16).Value = 250

# Row 16
$ws.Cells.Item(16, 4).Value = 44418
$ws.Cells.Item(16, 8).Value = "Española"
$ws.Cells.Item(16, 10).Value = 80
$ws.Cells.Item(16, 11).Value = 16000
$ws.Cells.Item(16, 12).Value = 16000
$ws.Cells.Item(16, 13).Value = 16000
$ws.Cells.Item(16, 14).Value = "$/caja 30 unidades"
$ws.Cells.Item(16, 16).Value = 533
$ws.Cells.Item(16, 17).Value = 30

# Row 17
$ws.Cells.Item(17, 4).Value = 44491
$ws.Cells.Item(17, 8).Value = "Madrigal"
$ws.Cells.Item(17, 10).Value = 200
$ws.Cells.Item(17, 11).Value = 11000
$ws.Cells.Item(17, 12).Value = 11000
$ws.Cells.Item(17, 13).Value = 11000
$ws.Cells.Item(17, 14).Value = "$/caja 40 unidades"
$ws.Cells.Item(17, 16).Value = 275
$ws.Cells.Item(17, 17).Value = 40
